$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price/volume refresh + FraxShare/PaxDollar row swap)
# D-column price values are prefixed with a leading apostrophe so Excel stores them
# as text (matching the source data, which uses "."-grouped strings, not numbers)
# rather than re-interpreting/rounding them as numeric values.

# Row 2
$ws.Range("D2").Value = '''26.079.25'
$ws.Range("E2").Value = '  -0.05%  '

# Row 3
$ws.Range("D3").Value = '''1.639.68'
$ws.Range("E3").Value = '  -1.60%  '

# Row 4
$ws.Range("E4").Value = '  -0.22%  '

# Row 5
$ws.Range("D5").Value = '''213.83'
$ws.Range("E5").Value = '  +2.15%  '

# Row 6
$ws.Range("E6").Value = '  -0.08%  '

# Row 7
$ws.Range("E7").Value = '  -0.18%  '

# Row 8
$ws.Range("D8").Value = '''0.2597'
$ws.Range("E8").Value = '  -0.88%  '

# Row 9
$ws.Range("D9").Value = '''0.06295'
$ws.Range("E9").Value = '  +0.21%  '

# Row 10
$ws.Range("D10").Value = '''20.65'
$ws.Range("E10").Value = '  -1.97%  '

# Row 11
$ws.Range("D11").Value = '''0.07656'
$ws.Range("E11").Value = '  +1.65%  '

# Row 12
$ws.Range("D12").Value = '''1.641.12'
$ws.Range("E12").Value = '  -1.59%  '

# Row 13
$ws.Range("E13").Value = '  -0.43%  '

# Row 14
$ws.Range("D14").Value = '''1.861.31'
$ws.Range("E14").Value = '  -1.71%  '

# Row 15
$ws.Range("D15").Value = '''0.5525'
$ws.Range("E15").Value = '  +0.49%  '

# Row 16
$ws.Range("D16").Value = '''0.0₅8253'
$ws.Range("E16").Value = '  +4.44%  '

# Row 17
$ws.Range("D17").Value = '''64.98'
$ws.Range("E17").Value = '  -2.16%  '

# Row 18
$ws.Range("D18").Value = '''26.066.49'
$ws.Range("E18").Value = '  -0.27%  '

# Row 19
$ws.Range("E19").Value = '  -0.15%  '

# Row 20
$ws.Range("D20").Value = '''4.686'
$ws.Range("E20").Value = '  -0.61%  '

# Row 21
$ws.Range("D21").Value = '''188.47'
$ws.Range("E21").Value = '  +1.25%  '

# Row 22
$ws.Range("D22").Value = '''10.17'
$ws.Range("E22").Value = '  -1.01%  '

# Row 23
$ws.Range("D23").Value = '''6.162'
$ws.Range("E23").Value = '  +0.12%  '

# Row 24
$ws.Range("D24").Value = '''1.001'
$ws.Range("E24").Value = '  -0.26%  '

# Row 25
$ws.Range("D25").Value = '''145.72'
$ws.Range("E25").Value = '  -2.47%  '

# Row 26
$ws.Range("D26").Value = '''0.1215'
$ws.Range("E26").Value = '  -2.40%  '

# Row 27
$ws.Range("D27").Value = '''7.417'
$ws.Range("E27").Value = '  -0.60%  '

# Row 28
$ws.Range("D28").Value = '''15.80'
$ws.Range("E28").Value = '  -0.47%  '

# Row 29
$ws.Range("D29").Value = '''1.397'
$ws.Range("E29").Value = '  +3.58%  '

# Row 30
$ws.Range("D30").Value = '''0.05972'
$ws.Range("E30").Value = '  -5.05%  '

# Row 31
$ws.Range("D31").Value = '''1.255'
$ws.Range("E31").Value = '  -1.36%  '

# Row 32
$ws.Range("D32").Value = '''3.439'
$ws.Range("E32").Value = '  -1.20%  '

# Row 33
$ws.Range("E33").Value = '  -0.04%  '

# Row 34
$ws.Range("D34").Value = '''1.642'
$ws.Range("E34").Value = '  +0.69%  '

# Row 35
$ws.Range("D35").Value = '''0.9839'
$ws.Range("E35").Value = '  -1.23%  '

# Row 36
$ws.Range("E36").Value = '  -0.54%  '

# Row 37
$ws.Range("E37").Value = '  +1.18%  '

# Row 38
$ws.Range("E38").Value = '  -5.58%  '

# Row 39
$ws.Range("D39").Value = '''0.01615'
$ws.Range("E39").Value = '  +0.19%  '

# Row 40
$ws.Range("D40").Value = '''0.8496'
$ws.Range("E40").Value = '  -2.50%  '

# Row 41
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '''1.001'
$ws.Range("E41").Value = '  -0.23%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''5.706'
$ws.Range("E42").Value = '  -6.15%  '

# Row 43
$ws.Range("D43").Value = '''1.034.58'
$ws.Range("E43").Value = '  -6.51%  '

# Row 44
$ws.Range("D44").Value = '''100.21'
$ws.Range("E44").Value = '  +0.47%  '

# Row 45
$ws.Range("D45").Value = '''1.787.67'
$ws.Range("E45").Value = '  -1.62%  '

# Row 46
$ws.Range("E46").Value = '  -2.78%  '

# Row 47
$ws.Range("D47").Value = '''55.80'
$ws.Range("E47").Value = '  +1.09%  '

# Row 48
$ws.Range("D48").Value = '''1.004'
$ws.Range("E48").Value = '  +0.38%  '

# Row 49
$ws.Range("D49").Value = '''8.049'
$ws.Range("E49").Value = '  +0.53%  '

# Row 50
$ws.Range("D50").Value = '''0.05157'
$ws.Range("E50").Value = '  -1.47%  '

# Row 51
$ws.Range("E51").Value = '  -0.70%  '
